# 1. Plane Support 2. Ready to item table 3. add critical
# Insert two new columns (Q,R) for "CriticalRate" / "CriticalDamage" ahead of
# the existing FireTime/ReloadTime/ShellCnt/maxDamage/Name/ResourceName block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at Q:R - shifts old Q..W (FireTime..blank) to S..Y
$ws.Columns("Q:R").Insert()

# New header cells
$ws.Range("Q1").Value = "CriticalRate"
$ws.Range("R1").Value = "CriticalDamage"

# New data values for the 10 populated data rows (rows 2-11)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 17).Value = 5
    $ws.Cells.Item($r, 18).Value = 5
}

# Column widths for the two new columns
$ws.Columns("Q").ColumnWidth = 9.66
$ws.Columns("R").ColumnWidth = 13.16

# Restore the selected cell to match the post-edit workbook state
$null = $ws.Range("Q21").Select()
